$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wallpapers")

$ws.Range("A9").Value = "birthday"
$ws.Range("B9").Value = "https://wallpaperscraft.com/download/pie_candles_birthday_1371/1080x1920"

$ws.Activate()
$ws.Range("B9").Select()
